$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 4649.4287
$ws.Range("J38").Value = 8084.1665
$ws.Range("L38").Value = 24252.4995
$ws.Range("N38").Value = -24996.4995
$ws.Range("H55").Value = 210.81818
$ws.Range("I55").Value = 112.8
$ws.Range("J55").Value = 292.5
$ws.Range("K55").Value = 112.8
$ws.Range("L55").Value = 292.5
$ws.Range("M55").Value = 101.2
$ws.Range("N55").Value = -720.5
$ws.Range("H98").Value = 1810.7858
$ws.Range("I98").Value = 1775.1
$ws.Range("K98").Value = 1775.1
$ws.Range("M98").Value = -277.0999999999999
$ws.Range("H122").Value = 1810.7858
$ws.Range("I122").Value = 1775.1
$ws.Range("K122").Value = 5325.299999999999
$ws.Range("M122").Value = -2875.299999999999
$ws.Range("H137").Value = 2150.4285
$ws.Range("I137").Value = 2126.575
$ws.Range("K137").Value = 6379.724999999999
$ws.Range("M137").Value = -3829.724999999999
$ws.Range("H138").Value = 7540.2
$ws.Range("J138").Value = 9488.416999999999
$ws.Range("L138").Value = 28465.251
$ws.Range("N138").Value = -38745.251

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2867.82
$ws.Range("I32").Value = 1718.6483
$ws.Range("J32").Value = 14487.223
$ws.Range("K32").Value = 1718.6483
$ws.Range("L32").Value = 14487.223
$ws.Range("M32").Value = -1431.6483
$ws.Range("N32").Value = -15061.223
$ws.Range("H45").Value = 3225.682
$ws.Range("I45").Value = 1697.6
$ws.Range("J45").Value = 4499.0835
$ws.Range("K45").Value = 1697.6
$ws.Range("L45").Value = 4499.0835
$ws.Range("M45").Value = -1320.6
$ws.Range("N45").Value = -5253.0835
$ws.Range("H74").Value = 8623909
$ws.Range("I74").Value = 10206583
$ws.Range("K74").Value = 10206583
$ws.Range("M74").Value = -10205709
$ws.Range("H77").Value = 8623909
$ws.Range("I77").Value = 10206583
$ws.Range("K77").Value = 51032915
$ws.Range("M77").Value = -51028547
$ws.Range("H132").Value = 3726.54
$ws.Range("I132").Value = 2178.0977
$ws.Range("K132").Value = 6534.293099999999
$ws.Range("M132").Value = -4004.293099999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1996.375
$ws.Range("I134").Value = 1480.129
$ws.Range("K134").Value = 4440.387
$ws.Range("M134").Value = -1905.387

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7288.7812
$ws.Range("I31").Value = 2984.3684
$ws.Range("J31").Value = 13579.846
$ws.Range("K31").Value = 2984.3684
$ws.Range("L31").Value = 13579.846
$ws.Range("M31").Value = -2689.3684
$ws.Range("N31").Value = -14169.846
$ws.Range("H34").Value = 7288.7812
$ws.Range("I34").Value = 2984.3684
$ws.Range("J34").Value = 13579.846
$ws.Range("K34").Value = 2984.3684
$ws.Range("L34").Value = 13579.846
$ws.Range("M34").Value = -2782.3684
$ws.Range("N34").Value = -13983.846
$ws.Range("H86").Value = 4622.364
$ws.Range("I86").Value = 4141.5
$ws.Range("J86").Value = 5199.4
$ws.Range("K86").Value = 4141.5
$ws.Range("L86").Value = 5199.4
$ws.Range("M86").Value = -3018.5
$ws.Range("N86").Value = -7445.4
$ws.Range("H89").Value = 4622.364
$ws.Range("I89").Value = 4141.5
$ws.Range("J89").Value = 5199.4
$ws.Range("K89").Value = 20707.5
$ws.Range("L89").Value = 25997
$ws.Range("M89").Value = -15091.5
$ws.Range("N89").Value = -37229
$ws.Range("H99").Value = 1404.7567
$ws.Range("I99").Value = 1374.4062
$ws.Range("K99").Value = 1374.4062
$ws.Range("M99").Value = 123.5938000000001
$ws.Range("H105").Value = 1988.5
$ws.Range("I105").Value = 2151.3333
$ws.Range("J105").Value = 1500
$ws.Range("K105").Value = 2151.3333
$ws.Range("L105").Value = 1500
$ws.Range("M105").Value = -404.3332999999998
$ws.Range("N105").Value = -4994
$ws.Range("H107").Value = 3008.0908
$ws.Range("I107").Value = 2886.375
$ws.Range("K107").Value = 2886.375
$ws.Range("M107").Value = -966.375
$ws.Range("H126").Value = 1404.7567
$ws.Range("I126").Value = 1374.4062
$ws.Range("K126").Value = 4123.2186
$ws.Range("M126").Value = -1653.2186

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2299.5715
$ws.Range("I5").Value = 3056.8572
$ws.Range("J5").Value = 1542.2858
$ws.Range("K5").Value = 9170.571599999999
$ws.Range("L5").Value = 4626.857400000001
$ws.Range("M5").Value = -9058.571599999999
$ws.Range("N5").Value = -4850.857400000001
$ws.Range("H11").Value = 760.7273
$ws.Range("I11").Value = 263.1111
$ws.Range("K11").Value = 789.3333
$ws.Range("M11").Value = -649.3333
$ws.Range("H12").Value = 224.90909
$ws.Range("J12").Value = 241.44444
$ws.Range("L12").Value = 724.33332
$ws.Range("N12").Value = -1070.33332
$ws.Range("H38").Value = 253.2
$ws.Range("I38").Value = 200.25
$ws.Range("K38").Value = 600.75
$ws.Range("M38").Value = -253.75
$ws.Range("H131").Value = 869372.5
$ws.Range("J131").Value = 1303581.1
$ws.Range("L131").Value = 3910743.3
$ws.Range("N131").Value = -3920823.3
$ws.Range("H135").Value = 2299.5715
$ws.Range("I135").Value = 3056.8572
$ws.Range("J135").Value = 1542.2858
$ws.Range("K135").Value = 27511.7148
$ws.Range("L135").Value = 13880.5722
$ws.Range("M135").Value = -24976.7148
$ws.Range("N135").Value = -18950.5722

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6850.1704
$ws.Range("I70").Value = 6204
$ws.Range("J70").Value = 6878.8887
$ws.Range("K70").Value = 6204
$ws.Range("L70").Value = 6878.8887
$ws.Range("M70").Value = -5934
$ws.Range("N70").Value = -7418.8887
$ws.Range("H73").Value = 6850.1704
$ws.Range("I73").Value = 6204
$ws.Range("J73").Value = 6878.8887
$ws.Range("K73").Value = 6204
$ws.Range("L73").Value = 6878.8887
$ws.Range("M73").Value = -5268
$ws.Range("N73").Value = -8750.8887
$ws.Range("H132").Value = 4792.948
$ws.Range("I132").Value = 4202.8696
$ws.Range("K132").Value = 12608.6088
$ws.Range("M132").Value = -10078.6088

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7568.0356
$ws.Range("I40").Value = 6621.0835
$ws.Range("J40").Value = 13249.75
$ws.Range("K40").Value = 6621.0835
$ws.Range("L40").Value = 13249.75
$ws.Range("M40").Value = -6485.0835
$ws.Range("N40").Value = -13521.75
$ws.Range("H122").Value = 3112.6487
$ws.Range("I122").Value = 3267.6086
$ws.Range("K122").Value = 9802.825800000001
$ws.Range("M122").Value = -7352.825800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3545.0476
$ws.Range("J96").Value = 4163.5454
$ws.Range("L96").Value = 4163.5454
$ws.Range("N96").Value = -6909.5454
